# Agenda: generacion de puntos.
# Rename the "Contenido con Imagen" slide layout to "Imagen".
$p = $ppt.ActivePresentation

$renamed = $false
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    if ($layout.Name -eq "Contenido con Imagen") {
        $layout.Name = "Imagen"
        $renamed = $true
    }
}

if (-not $renamed) {
    # Fallback: the "Contenido con Imagen" layout is known to be the 3rd
    # layout on the slide master (ppt/slideLayouts/slideLayout3.xml).
    $layouts.Item(3).Name = "Imagen"
}

Write-Host "Renamed layout 3 to:" $layouts.Item(3).Name
